$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-reg_center_machine")

# Update data values (2nd May Data Refresh)
$ws.Range("A3").Value = 10003
$ws.Range("A25").Value = 10003

# Update the selection/view: select entire row 34 (clears the old pinned top row / selection)
$ws.Activate()
$ws.Range("A34:XFD1048576").Select()
